# Commit: "Updated AB link to be VB link"
# The slide showing the system-model diagram has a small text label
# ("AB") identifying one of the UAV/ground-user links; rename it to "VB".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        if ($shape.TextFrame.TextRange.Text -eq "AB") {
            $target = $shape
            break
        }
    }
}

if ($target -eq $null) {
    # Fallback: known shape ("TextBox 1092") in case text already changed.
    $target = $s.Shapes.Item("TextBox 1092")
}

$target.TextFrame.TextRange.Text = "VB"
